$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the two unused/empty sheets (WorkActivity, Sheet4)
# ---------------------------------------------------------------------------
$wb.Worksheets("WorkActivity").Delete()
$wb.Worksheets("Sheet4").Delete()

# ---------------------------------------------------------------------------
# 2. Populate CompanyData with job history rows
#    (order of entry matters for shared-string indices: column A, then B,
#    then C only where C must be literal text - rows 9-11)
# ---------------------------------------------------------------------------
$cd = $wb.Worksheets("CompanyData")

$cd.Range("A2").Value = "CTO"
$cd.Range("B2").Value = "Hotpathz Inc"
$cd.Range("C2").Value = 94947
$cd.Range("D2").Value = 20180301

$cd.Range("A3").Value = "VP of Engineering"
$cd.Range("B3").Value = "Hotpathz Inc"
$cd.Range("C3").Value = 94947
$cd.Range("D3").Value = 20170901
$cd.Range("E3").Value = 20180301

$cd.Range("A4").Value = "Board Member Driver Rehabilitation Institute"
$cd.Range("C4").Value = 94947
$cd.Range("D4").Value = 20170901

$cd.Range("A5").Value = "President"
$cd.Range("B5").Value = "Black Shell Consulting"
$cd.Range("C5").Value = 76244
$cd.Range("D5").Value = 20151101

$cd.Range("A6").Value = "Business Analytics Analyst IV"
$cd.Range("B6").Value = "Beacon Hill Staffing"
$cd.Range("C6").Value = 75254
$cd.Range("D6").Value = 20161101
$cd.Range("E6").Value = 20170301

$cd.Range("A7").Value = "EMR/Soarian System Analyst"
$cd.Range("B7").Value = "Arkansas Heart Hospital"
$cd.Range("C7").Value = 72211
$cd.Range("D7").Value = 20130901
$cd.Range("E7").Value = 20151101

$cd.Range("A8").Value = "Consultant"
$cd.Range("B8").Value = "Health Data Specialists"
$cd.Range("C8").Value = 70070
$cd.Range("D8").Value = 20130901
$cd.Range("E8").Value = 20120201

# Remaining zip codes have leading zeros, so the column is switched to a text
# format before they're typed in (rows 9-11 land as literal text values).
$cd.Range("C1:C11").NumberFormat = "@"

$cd.Range("A9").Value = "Sr. Clinical Analyst"
$cd.Range("B9").Value = "Southern New Hampshire Medical Center"
$cd.Range("C9").Value = "03060"
$cd.Range("D9").Value = 20080301
$cd.Range("E9").Value = 20120601

$cd.Range("A10").Value = "Clinical Consultant"
$cd.Range("B10").Value = "Stoltenburg Consulting"
$cd.Range("C10").Value = "15106"
$cd.Range("D10").Value = 20100701
$cd.Range("E10").Value = 20120201

$cd.Range("A11").Value = "Implementation and Support Specialist"
$cd.Range("C11").Value = "06105"
$cd.Range("D11").Value = 20010101
$cd.Range("E11").Value = 20070601

# ---------------------------------------------------------------------------
# 3. Populate the "person" sheet with contact info + hyperlinked email
# ---------------------------------------------------------------------------
$person = $wb.Worksheets("person")
$person.Range("A2").Value = "Christopher"
$person.Range("B2").Value = "Huntley"
$person.Range("C2").Value = 94947
$person.Range("D2").Value = 4157746293
$person.Hyperlinks.Add($person.Range("E2"), "mailto:chuntley@blackshellconsulting.com")
$person.Range("E2").Value = "chuntley@blackshellconsulting.com"

# ---------------------------------------------------------------------------
# 4. New JobStart / JobEnd column headers on CompanyData
# ---------------------------------------------------------------------------
$cd.Range("D1").Value = "JobStart"
$cd.Range("E1").Value = "JobEnd"

# ---------------------------------------------------------------------------
# 5. Column widths on the "person" sheet
# ---------------------------------------------------------------------------
$person.Columns("A").ColumnWidth = 11.166666666666666
$person.Columns("B").ColumnWidth = 9.166666666666666
$person.Columns("D").ColumnWidth = 10.25

# ---------------------------------------------------------------------------
# 6. View state: selections per sheet, finishing with Skills as active tab
# ---------------------------------------------------------------------------
$person.Range("U38").Select()
$cd.Range("D12").Select()

$skills = $wb.Worksheets("Skills")
$skills.Range("B50").Select()
